$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.791.58'
$ws.Cells.Item(2, 5).Value = '  -0.54%  '
$ws.Cells.Item(3, 4).Value = '3.410.79'
$ws.Cells.Item(3, 5).Value = '  +0.74%  '
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '568.76'
$ws.Cells.Item(5, 5).Value = '  +0.52%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '157.34'
$ws.Cells.Item(6, 5).Value = '  +1.08%  '
$ws.Cells.Item(7, 5).Value = '  -0.08%  '
$ws.Cells.Item(8, 4).Value = '3.414.47'
$ws.Cells.Item(8, 5).Value = '  +0.78%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.571'
$ws.Cells.Item(9, 5).Value = '  -5.70%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.21'
$ws.Cells.Item(10, 5).Value = '  +0.75%  '
$ws.Cells.Item(11, 5).Value = '  -1.55%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.428'
$ws.Cells.Item(12, 5).Value = '  -2.74%  '
$ws.Cells.Item(13, 4).Value = '3.992.24'
$ws.Cells.Item(13, 5).Value = '  +0.35%  '
$ws.Cells.Item(14, 5).Value = '  -0.13%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '27.08'
$ws.Cells.Item(15, 5).Value = '  -1.72%  '
$ws.Cells.Item(16, 5).Value = '  -6.77%  '
$ws.Cells.Item(17, 4).Value = '63.871.18'
$ws.Cells.Item(17, 5).Value = '  -0.59%  '
$ws.Cells.Item(18, 4).Value = '3.406.31'
$ws.Cells.Item(18, 5).Value = '  +0.57%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.14'
$ws.Cells.Item(19, 5).Value = '  -2.77%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '13.62'
$ws.Cells.Item(20, 5).Value = '  -1.77%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '376.88'
$ws.Cells.Item(21, 5).Value = '  +0.93%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '7.73'
$ws.Cells.Item(22, 5).Value = '  -2.33%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  +0.25%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '71.31'
$ws.Cells.Item(24, 5).Value = '  -0.36%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.518'
$ws.Cells.Item(25, 5).Value = '  -5.11%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.0000116'
$ws.Cells.Item(26, 5).Value = '  -0.23%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '9.71'
$ws.Cells.Item(27, 5).Value = '  -1.43%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.178'
$ws.Cells.Item(28, 5).Value = '  +1.56%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.998'
$ws.Cells.Item(29, 5).Value = '  -0.22%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '6.05'
$ws.Cells.Item(30, 5).Value = '  +0.07%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.40'
$ws.Cells.Item(31, 5).Value = '  -3.42%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.98'
$ws.Cells.Item(32, 5).Value = '  -1.42%  '
$ws.Cells.Item(33, 5).Value = '  +0.08%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '22.84'
$ws.Cells.Item(34, 5).Value = '  -0.95%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '6.96'
$ws.Cells.Item(35, 5).Value = '  -2.16%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.54'
$ws.Cells.Item(36, 5).Value = '  -2.48%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '160.65'
$ws.Cells.Item(37, 5).Value = '  +0.44%  '
$ws.Cells.Item(38, 5).Value = '  -1.74%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.826'
$ws.Cells.Item(39, 5).Value = '  +8.52%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '26.10'
$ws.Cells.Item(40, 5).Value = '  -1.60%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0729'
$ws.Cells.Item(41, 5).Value = '  -3.51%  '
$ws.Cells.Item(42, 4).Value = '2.800.04'
$ws.Cells.Item(42, 5).Value = '  -0.83%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '42.69'
$ws.Cells.Item(43, 5).Value = '  +0.29%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '4.41'
$ws.Cells.Item(44, 5).Value = '  -3.92%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '6.38'
$ws.Cells.Item(45, 5).Value = '  -5.09%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '25.75'
$ws.Cells.Item(46, 5).Value = '  +1.42%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0305'
$ws.Cells.Item(47, 5).Value = '  -2.21%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.35'
$ws.Cells.Item(48, 5).Value = '  +10.85%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '323.84'
$ws.Cells.Item(49, 5).Value = '  +4.55%  '
$ws.Cells.Item(50, 5).Value = '  -2.70%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '6.32'
$ws.Cells.Item(51, 5).Value = '  -2.93%  '
